# Apply updated cryptocurrency price/volume data to the worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.040.68"
$ws.Range("E2").Value = "  +0.41%  "
$ws.Range("D3").Value = "2.663.85"
$ws.Range("E3").Value = "  +1.72%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "530.10"
$ws.Range("E5").Value = "  +3.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "155.62"
$ws.Range("E6").Value = "  +0.23%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.581"
$ws.Range("E8").Value = "  -1.21%  "
$ws.Range("E9").Value = "  -4.40%  "
$ws.Range("E10").Value = "  +3.73%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.352"
$ws.Range("E11").Value = "  +1.49%  "
$ws.Range("E12").Value = "  -0.48%  "
$ws.Range("D13").Value = "3.131.64"
$ws.Range("E13").Value = "  +1.83%  "
$ws.Range("D14").Value = "60.977.09"
$ws.Range("E14").Value = "  +0.41%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "22.09"
$ws.Range("E15").Value = "  +1.71%  "
$ws.Range("E16").Value = "  +1.05%  "
$ws.Range("D17").Value = "2.672.01"
$ws.Range("E17").Value = "  +1.76%  "
$ws.Range("E18").Value = "  +0.90%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "354.35"
$ws.Range("E19").Value = "  -0.82%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.69"
$ws.Range("E20").Value = "  +0.66%  "
$ws.Range("E21").Value = "  +2.22%  "
$ws.Range("E22").Value = "  +0.23%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "61.82"
$ws.Range("E23").Value = "  +1.78%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.432"
$ws.Range("E24").Value = "  +1.78%  "
$ws.Range("E25").Value = "  +0.85%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  +0.26%  "
$ws.Range("D27").Value = "0.0₃0859"
$ws.Range("E27").Value = "  +1.22%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.34"
$ws.Range("E28").Value = "  -0.59%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.19"
$ws.Range("E30").Value = "  +4.29%  "
$ws.Range("E31").Value = "  +0.38%  "
$ws.Range("E32").Value = "  +2.69%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "149.65"
$ws.Range("E33").Value = "  -1.71%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.12"
$ws.Range("E34").Value = "  +2.71%  "
$ws.Range("E35").Value = "  +0.34%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.920"
$ws.Range("E36").Value = "  +8.06%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.892"
$ws.Range("E37").Value = "  +1.48%  "
$ws.Range("B38").Value = "Bittensor"
$ws.Range("C38").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "307.57"
$ws.Range("E38").Value = "  +4.10%  "
$ws.Range("B39").Value = "OKB"
$ws.Range("C39").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.83"
$ws.Range("E39").Value = "  +1.38%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.49"
$ws.Range("E40").Value = "  -0.39%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.78"
$ws.Range("E41").Value = "  +0.47%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.648"
$ws.Range("E42").Value = "  +3.79%  "
$ws.Range("B43").Value = "Stellar"
$ws.Range("C43").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.102"
$ws.Range("E43").Value = "  +0.02%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "20.40"
$ws.Range("E44").Value = "  +2.42%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.996"
$ws.Range("E46").Value = "  +0.08%  "
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0242"
$ws.Range("E47").Value = "  +2.82%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.90"
$ws.Range("E48").Value = "  -0.99%  "
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "19.18"
$ws.Range("E49").Value = "  +6.14%  "
$ws.Range("B50").Value = "WhiteBITCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "10.35"
$ws.Range("E50").Value = "  +0.38%  "
$ws.Range("D51").Value = "2.001.23"
$ws.Range("E51").Value = "  -0.02%  "
